# Auto-generated edit script: apply numeric updates from the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 50236
$ws.Range("J87").Value = 50236
$ws.Range("L87").Value = 50236
$ws.Range("N87").Value = -52732

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 50236
$ws.Range("J90").Value = 50236
$ws.Range("L90").Value = 150708
$ws.Range("N90").Value = -163188

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 117601.3
$ws.Range("I132").Value = 2467.1323
$ws.Range("J132").Value = 529660.4399999999
$ws.Range("K132").Value = 7401.396900000001
$ws.Range("L132").Value = 1588981.32
$ws.Range("M132").Value = -4871.396900000001
$ws.Range("N132").Value = -1594041.32

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 29209.281
$ws.Range("I137").Value = 42618
$ws.Range("J137").Value = 7755.3335
$ws.Range("K137").Value = 127854
$ws.Range("L137").Value = 23266.0005
$ws.Range("M137").Value = -125304
$ws.Range("N137").Value = -28366.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1418.0273
$ws.Range("I141").Value = 855.0179000000001
$ws.Range("J141").Value = 3272.647
$ws.Range("K141").Value = 2565.0537
$ws.Range("L141").Value = 9817.940999999999
$ws.Range("M141").Value = 2614.9463
$ws.Range("N141").Value = -20177.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5331.9
$ws.Range("I32").Value = 5409.284
$ws.Range("J32").Value = 4635.4443
$ws.Range("K32").Value = 5409.284
$ws.Range("L32").Value = 4635.4443
$ws.Range("M32").Value = -5122.284
$ws.Range("N32").Value = -5209.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1347.9153
$ws.Range("I61").Value = 1451.4166
$ws.Range("J61").Value = 896.2727
$ws.Range("K61").Value = 1451.4166
$ws.Range("L61").Value = 896.2727
$ws.Range("M61").Value = -1239.4166
$ws.Range("N61").Value = -1320.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 18391.861
$ws.Range("I74").Value = 21443.959
$ws.Range("J74").Value = 1774.8889
$ws.Range("K74").Value = 21443.959
$ws.Range("L74").Value = 1774.8889
$ws.Range("M74").Value = -20569.959
$ws.Range("N74").Value = -3522.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 18391.861
$ws.Range("I77").Value = 21443.959
$ws.Range("J77").Value = 1774.8889
$ws.Range("K77").Value = 107219.795
$ws.Range("L77").Value = 8874.4445
$ws.Range("M77").Value = -102851.795
$ws.Range("N77").Value = -17610.4445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1347.9153
$ws.Range("I136").Value = 1451.4166
$ws.Range("J136").Value = 896.2727
$ws.Range("K136").Value = 4354.2498
$ws.Range("L136").Value = 2688.8181
$ws.Range("M136").Value = -1804.2498
$ws.Range("N136").Value = -7788.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 16597.408
$ws.Range("I134").Value = 929.24
$ws.Range("J134").Value = 53902.57
$ws.Range("K134").Value = 2787.72
$ws.Range("L134").Value = 161707.71
$ws.Range("M134").Value = -252.7200000000003
$ws.Range("N134").Value = -166777.71

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6660.9316
$ws.Range("I31").Value = 4992.356
$ws.Range("J31").Value = 13692.786
$ws.Range("K31").Value = 4992.356
$ws.Range("L31").Value = 13692.786
$ws.Range("M31").Value = -4697.356
$ws.Range("N31").Value = -14282.786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6660.9316
$ws.Range("I34").Value = 4992.356
$ws.Range("J34").Value = 13692.786
$ws.Range("K34").Value = 4992.356
$ws.Range("L34").Value = 13692.786
$ws.Range("M34").Value = -4790.356
$ws.Range("N34").Value = -14096.786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 854.2537
$ws.Range("I58").Value = 657.4039
$ws.Range("J58").Value = 1536.6666
$ws.Range("K58").Value = 657.4039
$ws.Range("L58").Value = 1536.6666
$ws.Range("M58").Value = -454.4039
$ws.Range("N58").Value = -1942.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 776.5893
$ws.Range("I132").Value = 828.2
$ws.Range("J132").Value = 346.5
$ws.Range("K132").Value = 2484.6
$ws.Range("L132").Value = 1039.5
$ws.Range("M132").Value = 45.39999999999964
$ws.Range("N132").Value = -6099.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 749.9355
$ws.Range("I134").Value = 713.9778
$ws.Range("J134").Value = 845.1177
$ws.Range("K134").Value = 2141.9334
$ws.Range("L134").Value = 2535.3531
$ws.Range("M134").Value = 393.0666000000001
$ws.Range("N134").Value = -7605.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 854.2537
$ws.Range("I136").Value = 657.4039
$ws.Range("J136").Value = 1536.6666
$ws.Range("K136").Value = 1972.2117
$ws.Range("L136").Value = 4609.9998
$ws.Range("M136").Value = 577.7882999999999
$ws.Range("N136").Value = -9709.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 2000
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 50586696
$ws.Range("I112").Value = 1990
$ws.Range("J112").Value = 67448264
$ws.Range("K112").Value = 5970
$ws.Range("L112").Value = 202344792
$ws.Range("M112").Value = -4862
$ws.Range("N112").Value = -202347008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3600
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 3600
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 10800
$ws.Range("N115").Value = -13150
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 12963
$ws.Range("I120").Value = 4205.8
$ws.Range("J120").Value = 15395.556
$ws.Range("K120").Value = 12617.4
$ws.Range("L120").Value = 46186.66800000001
$ws.Range("M120").Value = -7779.400000000001
$ws.Range("N120").Value = -55862.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 96154550
$ws.Range("J131").Value = 125000776
$ws.Range("L131").Value = 375002328
$ws.Range("N131").Value = -375012408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 23729.021
$ws.Range("I132").Value = 1413.5186
$ws.Range("J132").Value = 57202.277
$ws.Range("K132").Value = 4240.5558
$ws.Range("L132").Value = 171606.831
$ws.Range("M132").Value = -1710.5558
$ws.Range("N132").Value = -176666.831

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 238355.7
$ws.Range("I132").Value = 54560.316
$ws.Range("J132").Value = 775603.75
$ws.Range("K132").Value = 163680.948
$ws.Range("L132").Value = 2326811.25
$ws.Range("M132").Value = -161150.948
$ws.Range("N132").Value = -2331871.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 936703.1
$ws.Range("I136").Value = 1253584.5
$ws.Range("K136").Value = 3760753.5
$ws.Range("M136").Value = -3758203.5
